$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "247.75"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.34"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.236"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05688"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.419"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.310"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8065"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8657"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1417"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07430"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03052"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03078"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09396"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.887"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001582"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04795"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.01828"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0005812"
$ws.Range("E19").Value = "18OneONEWorstin24h"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006372"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.005041"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0009969"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0001501"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.186"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3245"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1341"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03960"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006744"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1065"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003202"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008764"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005585"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4502"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1451"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
